# Auto-generated Excel COM-interop script to apply the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-E hold text data (coin name, link, price, % volume). Coerce to Text first so
# numeric-looking strings (e.g. "9.38") are NOT auto-converted to numbers by COM's usual
# type inference, matching the source workbook where every data cell is an inline string.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "61.927.72"
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
$ws.Cells.Item(3, 4).Value = "3.425.89"
$ws.Cells.Item(3, 5).Value = "  -0.22%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).Value = "410.55"
$ws.Cells.Item(5, 5).Value = "  +0.36%  "
$ws.Cells.Item(6, 4).Value = "129.89"
$ws.Cells.Item(6, 5).Value = "  +0.32%  "
$ws.Cells.Item(7, 4).Value = "0.634"
$ws.Cells.Item(7, 5).Value = "  -0.34%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 4).Value = "0.736"
$ws.Cells.Item(9, 5).Value = "  -2.45%  "
$ws.Cells.Item(10, 5).Value = "  -1.69%  "
$ws.Cells.Item(11, 4).Value = "43.60"
$ws.Cells.Item(11, 5).Value = "  +1.40%  "
$ws.Cells.Item(12, 4).Value = "0.0000226"
$ws.Cells.Item(12, 5).Value = "  +18.35%  "
$ws.Cells.Item(13, 4).Value = "9.38"
$ws.Cells.Item(13, 5).Value = "  +6.15%  "
$ws.Cells.Item(14, 4).Value = "3.964.80"
$ws.Cells.Item(14, 5).Value = "  -0.28%  "
$ws.Cells.Item(15, 5).Value = "  +0.36%  "
$ws.Cells.Item(16, 4).Value = "21.23"
$ws.Cells.Item(16, 5).Value = "  +3.95%  "
$ws.Cells.Item(17, 4).Value = "3.434.83"
$ws.Cells.Item(17, 5).Value = "  -1.22%  "
$ws.Cells.Item(18, 2).Value = "Polygon"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(18, 4).Value = "1.09"
$ws.Cells.Item(18, 5).Value = "  +3.17%  "
$ws.Cells.Item(19, 2).Value = "Uniswap"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(19, 4).Value = "12.36"
$ws.Cells.Item(19, 5).Value = "  +8.72%  "
$ws.Cells.Item(20, 4).Value = "61.865.04"
$ws.Cells.Item(20, 5).Value = "  -0.41%  "
$ws.Cells.Item(21, 4).Value = "513.86"
$ws.Cells.Item(21, 5).Value = "  +31.57%  "
$ws.Cells.Item(22, 4).Value = "92.33"
$ws.Cells.Item(22, 5).Value = "  +3.95%  "
$ws.Cells.Item(23, 4).Value = "3.34"
$ws.Cells.Item(23, 5).Value = "  +4.72%  "
$ws.Cells.Item(24, 4).Value = "13.47"
$ws.Cells.Item(24, 5).Value = "  +0.72%  "
$ws.Cells.Item(25, 5).Value = "  +3.55%  "
$ws.Cells.Item(26, 4).Value = "34.87"
$ws.Cells.Item(26, 5).Value = "  +8.50%  "
$ws.Cells.Item(27, 4).Value = "9.26"
$ws.Cells.Item(27, 5).Value = "  +8.74%  "
$ws.Cells.Item(28, 4).Value = "7.63"
$ws.Cells.Item(28, 5).Value = "  -0.79%  "
$ws.Cells.Item(29, 4).Value = "12.19"
$ws.Cells.Item(29, 5).Value = "  +3.36%  "
$ws.Cells.Item(30, 4).Value = "2.71"
$ws.Cells.Item(30, 5).Value = "  -0.67%  "
$ws.Cells.Item(31, 4).Value = "0.115"
$ws.Cells.Item(31, 5).Value = "  -1.10%  "
$ws.Cells.Item(32, 5).Value = "  -2.00%  "
$ws.Cells.Item(33, 4).Value = "41.99"
$ws.Cells.Item(33, 5).Value = "  -4.61%  "
$ws.Cells.Item(34, 4).Value = "59.74"
$ws.Cells.Item(34, 5).Value = "  +14.19%  "
$ws.Cells.Item(35, 5).Value = "  +0.03%  "
$ws.Cells.Item(36, 4).Value = "0.0500"
$ws.Cells.Item(36, 5).Value = "  +1.94%  "
$ws.Cells.Item(37, 5).Value = "  -0.04%  "
$ws.Cells.Item(38, 2).Value = "Stellar"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(38, 4).Value = "0.140"
$ws.Cells.Item(38, 5).Value = "  +5.01%  "
$ws.Cells.Item(39, 2).Value = "LidoDAOToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(39, 4).Value = "3.48"
$ws.Cells.Item(39, 5).Value = "  +3.51%  "
$ws.Cells.Item(40, 4).Value = "2.76"
$ws.Cells.Item(40, 5).Value = "  +19.19%  "
$ws.Cells.Item(41, 4).Value = "147.76"
$ws.Cells.Item(41, 5).Value = "  +4.53%  "
$ws.Cells.Item(42, 2).Value = "ARBITRUM"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(42, 4).Value = "2.11"
$ws.Cells.Item(42, 5).Value = "  +7.49%  "
$ws.Cells.Item(43, 4).Value = "2.93"
$ws.Cells.Item(43, 5).Value = "  +0.43%  "
$ws.Cells.Item(44, 2).Value = "TheGraph"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(44, 4).Value = "0.319"
$ws.Cells.Item(44, 5).Value = "  +2.69%  "
$ws.Cells.Item(45, 5).Value = "  +8.59%  "
$ws.Cells.Item(46, 4).Value = "16.78"
$ws.Cells.Item(46, 5).Value = "  +0.54%  "
$ws.Cells.Item(47, 4).Value = "2.34"
$ws.Cells.Item(47, 5).Value = "  +20.93%  "
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).Value = "23.20"
$ws.Cells.Item(48, 5).Value = "  +5.65%  "
$ws.Cells.Item(49, 2).Value = "BitcoinSV"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(49, 4).Value = "120.58"
$ws.Cells.Item(49, 5).Value = "  +26.66%  "
$ws.Cells.Item(50, 4).Value = "0.145"
$ws.Cells.Item(50, 5).Value = "  +18.73%  "
$ws.Cells.Item(51, 4).Value = "2.141.69"
$ws.Cells.Item(51, 5).Value = "  +1.12%  "

# Restore the default (General) style on the touched range so no stray number-format
# is left behind on the cells (mirrors the original "no explicit style" cells).
$dataRange.Style = "Normal"
